$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.929.42"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.352.40"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "240.57"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "0.661"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("D7").Value = "74.06"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "59.75"
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").Value = "33.27"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.109"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "7.28"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "2.703.11"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "16.21"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "0.905"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "2.357.07"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "43.862.20"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "78.20"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").Value = "253.11"
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("D25").Value = "3.80"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "2.51"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "10.47"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +3.05%  "
$ws.Range("D30").Value = "176.49"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "22.31"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "0.0751"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.09"
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "5.39"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "3.82"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").Value = "2.40"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0273"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").Value = "5.45"
$ws.Range("E41").Value = "  +14.39%  "
$ws.Range("D42").Value = "64.79"
$ws.Range("E42").Value = "  +15.55%  "
$ws.Range("D43").Value = "9.16"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").Value = "19.08"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.106"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.200"
$ws.Range("E46").Value = "  -5.88%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").Value = "98.55"
$ws.Range("E51").Value = "  -2.06%  "
